$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell I8: "colors", bold Calibri font, centered
$hdr = $ws.Range("I8")
$hdr.Value = "colors"
$hdr.Style = "Normal"
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108

# Boolean values for I9:I127 (TRUE/FALSE "colors" flag per gene row)
$vals = @($false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$true,$true,$false,$false,$false,$false,$true,$true,$true,$false,$true,$false,$false,$false,$true,$true,$true,$false,$false,$false,$false,$false,$false,$false,$true,$true,$true,$true,$true,$false,$false,$false,$true,$true,$false,$true,$false,$false,$false,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$true,$true,$true,$false,$false,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$true,$false,$true,$false,$false,$true,$false,$false,$false,$false,$true,$true,$false,$false,$true,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$false,$false,$true,$true,$false,$true,$true)

for ($i = 0; $i -lt $vals.Length; $i++) {
    $row = 9 + $i
    $cell = $ws.Range("I$row")
    $cell.Value = $vals[$i]
    $cell.Style = "Normal"
}

# Reflect the author's final selection/cursor position
[void]$ws.Range("K11").Select()
